# Primer commit rama Testint despues de modificar los archivos
#
# Original paragraph 1 text "Hola mundo" is re-split into proofed runs
# (spell-checker w:proofErr spellStart/spellEnd bookmarks around each
# word, keeping the visible text identical), a new blank paragraph is
# added, and a new paragraph "Hola, como va todo" is appended (also
# proof-marked), carrying forward the _GoBack bookmark to its end.
#
# w:proofErr is not reachable through a dedicated Word object-model
# property, so we drive it with Range.InsertXML, which accepts a raw
# WordprocessingML fragment (the same "single file XML"/WordOpenXML
# package format Range.WordOpenXML returns) and splices it in as real
# document content -- proofErr elements included -- at the target Range.

$d = $word.ActiveDocument

$wordOpenXmlHeader = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$wordOpenXmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rPr = '<w:rPr><w:lang w:val="en-US"/></w:rPr>'

# Paragraph 1: "Hola mundo" -> proof-marked "Hola" + " " + "mundo"
$para1 = '<w:p w:rsidR="00C20F83" w:rsidRPr="00744473" w:rsidRDefault="00744473">' +
  '<w:pPr>' + $rPr + '</w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r>' + $rPr + '<w:t>Hola</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r>' + $rPr + '<w:t>mundo</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'

# Paragraph 2: new blank paragraph
$para2 = '<w:p><w:pPr>' + $rPr + '</w:pPr></w:p>'

# Paragraph 3: new "Hola, como va todo" paragraph, carrying the
# _GoBack bookmark that used to sit at the end of paragraph 1.
$para3 = '<w:p>' +
  '<w:pPr>' + $rPr + '</w:pPr>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r>' + $rPr + '<w:t>Hola</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r>' + $rPr + '<w:t xml:space="preserve">, </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r>' + $rPr + '<w:t>como</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r>' + $rPr + '<w:t>va</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> todo</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>'

$xml = $wordOpenXmlHeader + $para1 + $para2 + $para3 + $wordOpenXmlFooter

# Replace the whole (only) paragraph's range -- this both rewrites its
# runs with the proof-marked split and appends the two new paragraphs
# after it, moving the _GoBack bookmark to the new last paragraph.
$r = $d.Paragraphs(1).Range
$r.InsertXML($xml)
